$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while forcing text storage
# (prevents Excel from re-interpreting numeric-looking strings, e.g. "0.520",
# as actual numbers, which would lose the trailing zero / exact formatting).
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "27.041.47"
$ws.Range("E2").Value = "  +2.78%  "

$ws.Range("D3").Value = "1.653.43"
$ws.Range("E3").Value = "  +3.62%  "

$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "215.25"
$ws.Range("E5").Value = "  +1.67%  "

$ws.Range("E6").Value = "  +1.67%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  +1.73%  "

$ws.Range("E9").Value = "  +1.66%  "

Set-TextValue $ws.Range("D10") "19.93"
$ws.Range("E10").Value = "  +4.68%  "

$ws.Range("E11").Value = "  +1.32%  "

$ws.Range("D12").Value = "1.887.21"
$ws.Range("E12").Value = "  +3.67%  "

$ws.Range("D13").Value = "1.648.81"
$ws.Range("E13").Value = "  +2.51%  "

$ws.Range("E14").Value = "  +2.30%  "

Set-TextValue $ws.Range("D15") "0.520"
$ws.Range("E15").Value = "  +3.21%  "

Set-TextValue $ws.Range("D16") "65.37"

Set-TextValue $ws.Range("D17") "240.08"
$ws.Range("E17").Value = "  +4.53%  "

$ws.Range("D18").Value = "27.039.50"
$ws.Range("E18").Value = "  +2.78%  "

Set-TextValue $ws.Range("D19") "7.84"
$ws.Range("E19").Value = "  +2.66%  "

$ws.Range("E20").Value = "  +1.33%  "

Set-TextValue $ws.Range("D21") "0.999"

$ws.Range("E22").Value = "  +4.46%  "

$ws.Range("E23").Value = "  +2.94%  "

$ws.Range("E24").Value = "  +3.84%  "

Set-TextValue $ws.Range("D25") "146.12"
$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("E27").Value = "  +2.26%  "

$ws.Range("E28").Value = "  +1.85%  "

$ws.Range("E30").Value = "  +0.84%  "

$ws.Range("E31").Value = "  +1.94%  "

$ws.Range("E32").Value = "  +3.36%  "

$ws.Range("D33").Value = "1.521.28"
$ws.Range("E33").Value = "  +1.56%  "

$ws.Range("E34").Value = "  +5.49%  "

Set-TextValue $ws.Range("D35") "1.60"
$ws.Range("E35").Value = "  +8.76%  "

$ws.Range("E36").Value = "  -0.38%  "

$ws.Range("E37").Value = "  +2.40%  "

Set-TextValue $ws.Range("D38") "0.889"

$ws.Range("E39").Value = "  +3.43%  "

$ws.Range("E40").Value = "  +3.00%  "

$ws.Range("E42").Value = "  +4.14%  "

Set-TextValue $ws.Range("D43") "65.85"
$ws.Range("E43").Value = "  +8.65%  "

$ws.Range("D44").Value = "1.793.31"
$ws.Range("E44").Value = "  +3.42%  "

$ws.Range("E45").Value = "  +2.06%  "

Set-TextValue $ws.Range("D46") "0.915"
$ws.Range("E46").Value = "  -2.32%  "

Set-TextValue $ws.Range("D47") "89.64"
$ws.Range("E47").Value = "  +1.32%  "

$ws.Range("E48").Value = "  +0.65%  "

$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("E51").Value = "  +2.19%  "
